# Weekly fruit/vegetable price update for "Femacal de La Calera - Esparragos".
# Reassigns date (D), quality (I), volume/price (J,K,L,M,P) and origin (O)
# values across the existing data rows (2-28) to their corrected weekly figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44172
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 1300
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = 1300
$ws.Range("P2").Value = 1300

# Row 3
$ws.Range("D3").Value = 44172
$ws.Range("J3").Value = 550

# Row 4
$ws.Range("D4").Value = 44159
$ws.Range("J4").Value = 1100
$ws.Range("K4").Value = 1300
$ws.Range("M4").Value = 1300
$ws.Range("P4").Value = 1300

# Row 5
$ws.Range("D5").Value = 44159
$ws.Range("J5").Value = 800

# Row 6
$ws.Range("D6").Value = 44167
$ws.Range("J6").Value = 1430
$ws.Range("K6").Value = 1200
$ws.Range("M6").Value = 1248
$ws.Range("P6").Value = 1248

# Row 7
$ws.Range("D7").Value = 44167
$ws.Range("J7").Value = 350

# Row 8
$ws.Range("D8").Value = 44169
$ws.Range("J8").Value = 950

# Row 9
$ws.Range("D9").Value = 44169
$ws.Range("J9").Value = 800

# Row 12
$ws.Range("D12").Value = 44174
$ws.Range("J12").Value = 2800
$ws.Range("K12").Value = 1200
$ws.Range("L12").Value = 1250
$ws.Range("M12").Value = 1221
$ws.Range("P12").Value = 1221

# Row 13
$ws.Range("D13").Value = 44174
$ws.Range("J13").Value = 1300

# Row 14
$ws.Range("D14").Value = 44168
$ws.Range("J14").Value = 1200

# Row 15
$ws.Range("D15").Value = 44168
$ws.Range("J15").Value = 850

# Row 16
$ws.Range("D16").Value = 44175
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1300
$ws.Range("L16").Value = 1300
$ws.Range("M16").Value = 1300
$ws.Range("P16").Value = 1300

# Row 17
$ws.Range("D17").Value = 44175
$ws.Range("J17").Value = 1450

# Row 18
$ws.Range("D18").Value = 44161
$ws.Range("J18").Value = 1600

# Row 19
$ws.Range("D19").Value = 44161
$ws.Range("J19").Value = 1850
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = 1000
$ws.Range("P19").Value = 1000

# Row 20
$ws.Range("D20").Value = 44165
$ws.Range("J20").Value = 720
$ws.Range("O20").Value = "Provincia de Quillota"

# Row 21
$ws.Range("D21").Value = 44165
$ws.Range("I21").Value = "Segunda"
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 1000
$ws.Range("P21").Value = 1000

# Row 22
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 1300
$ws.Range("L22").Value = 1300
$ws.Range("M22").Value = 1300
$ws.Range("P22").Value = 1300

# Row 23
$ws.Range("D23").Value = 44160
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 850
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = 1000
$ws.Range("P23").Value = 1000

# Row 24
$ws.Range("D24").Value = 44181
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 1000
$ws.Range("K24").Value = 1300
$ws.Range("L24").Value = 1300
$ws.Range("M24").Value = 1300
$ws.Range("P24").Value = 1300

# Row 25
$ws.Range("D25").Value = 44181
$ws.Range("I25").Value = "Segunda"
$ws.Range("J25").Value = 900
$ws.Range("K25").Value = 900
$ws.Range("L25").Value = 900
$ws.Range("M25").Value = 900
$ws.Range("P25").Value = 900

# Row 26
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 1200
$ws.Range("K26").Value = 1300
$ws.Range("L26").Value = 1300
$ws.Range("M26").Value = 1300
$ws.Range("P26").Value = 1300

# Row 27
$ws.Range("D27").Value = 44162
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = 1000
$ws.Range("P27").Value = 1000

# Row 28
$ws.Range("D28").Value = 44179
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 980
$ws.Range("K28").Value = 1200
$ws.Range("L28").Value = 1200
$ws.Range("M28").Value = 1200
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 1200
